# Tasks day 10 - actions and cookies
# Add a new "TestAutomation" worksheet (after the existing sheets) with a
# simple header row: Name | Price | location

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so it lands at the end
# of the tab strip (Sayfa1, Sayfa2, TestAutomation).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "TestAutomation"

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Price"
$ws.Range("C1").Value = "location"

# Keep the originally-active sheet (Sayfa2) selected, matching the source
# workbook's activeTab.
$wb.Worksheets.Item("Sayfa2").Activate()
